# Assign task owners ("bu" / "fu") to the plan sheet and highlight the
# header row of the reference "list" sheet; also restore the last-used
# selections on both sheets.

$wb = $excel.ActiveWorkbook
$planWs = $wb.Worksheets.Item("plan")
$listWs = $wb.Worksheets.Item("list")

# ---------------------------------------------------------------------
# 1. Fill in the "担当" (person in charge) column on the plan sheet.
#    Rows 5-11 -> "bu", rows 12-21 -> "fu"
# ---------------------------------------------------------------------
$buRows = 5..11
foreach ($r in $buRows) {
    $planWs.Range("J$r").Value = "bu"
}

$fuRows = 12..21
foreach ($r in $fuRows) {
    $planWs.Range("J$r").Value = "fu"
}

# ---------------------------------------------------------------------
# 2. Highlight the header row of the "list" sheet (A1:D1) with the
#    "Blue, Accent 5, Lighter 80%" fill color.
# ---------------------------------------------------------------------
$listWs.Range("A1:D1").Interior.Color = 15983578

# ---------------------------------------------------------------------
# 3. Restore the active selections that were left on each sheet.
# ---------------------------------------------------------------------
$listWs.Activate() | Out-Null
$listWs.Range("C13").Select() | Out-Null

$planWs.Activate() | Out-Null
$planWs.Range("E21:I21").Select() | Out-Null
